$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be stored as text so that numeric-looking
# strings (e.g. "25.08", "101.10") are not silently coerced to numbers and
# do not lose formatting such as trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.350.86'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.841.16'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '239.85'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '0.6289'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.07456'
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = '25.08'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.2897'
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = '1.833.39'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '4.973'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '0.00001034'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '81.84'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '6.238'
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '29.353.47'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '228.90'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").Value = '0.9996'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '7.378'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '158.02'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Value = '8.533'
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").Value = '0.1348'
$ws.Range("E26").Value = '  -1.78%  '
$ws.Range("D27").Value = '17.43'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '0.06874'
$ws.Range("E28").Value = '  +6.39%  '
$ws.Range("D29").Value = '1.454'
$ws.Range("E29").Value = '  +4.31%  '
$ws.Range("D30").Value = '1.490'
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").Value = '4.055'
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").Value = '1.829'
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '0.6986'
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '0.01844'
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '1.238.13'
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").Value = '6.804'
$ws.Range("E40").Value = '  +4.44%  '
$ws.Range("D41").Value = '0.9399'
$ws.Range("E41").Value = '  +3.37%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '101.10'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '1.980.21'
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("D45").Value = '65.40'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '7.038'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = '1.711'
$ws.Range("E48").Value = '  +2.30%  '
$ws.Range("D49").Value = '8.956'
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("D51").Value = '0.3910'
$ws.Range("E51").Value = '  -0.83%  '

# Restore the default General format/style so no stray style diffs are left behind.
$ws.Range("D2:E51").NumberFormat = "General"
$ws.Range("D2:E51").Style = "Normal"

